# Apply row 2-32 "name_duty" assignments (values + fonts) to match the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row is unchanged in content ("day" / "name_duty"); touch it defensively too.
$ws.Range("A1").Value = "day"
$ws.Range("B1").Value = "name_duty"

$ws.Range("B2").Value = 'なし'
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 10

$ws.Range("B3").Value = 'なし'
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Size = 10

$ws.Range("B4").Value = 'なし'
$ws.Range("B4").Font.Name = "Arial"
$ws.Range("B4").Font.Size = 10
$ws.Range("B4").Font.Color = 0

$ws.Range("B5").Value = 'なし'
$ws.Range("B5").Font.Name = "Arial"
$ws.Range("B5").Font.Size = 10
$ws.Range("B5").Font.Color = 0

$ws.Range("B6").Value = '白岩詩佑介'
$ws.Range("B6").Font.Name = "Arial"
$ws.Range("B6").Font.Size = 10

$ws.Range("B7").Value = '石井海成'
$ws.Range("B7").Font.Name = "Arial"
$ws.Range("B7").Font.Size = 10

$ws.Range("B8").Value = '林健太郎'
$ws.Range("B8").Font.Name = "Arial"
$ws.Range("B8").Font.Size = 10

$ws.Range("B9").Value = '小溝賢'
$ws.Range("B9").Font.Name = "Arial"
$ws.Range("B9").Font.Size = 10

$ws.Range("B10").Value = '小野文哉'
$ws.Range("B10").Font.Name = "Arial"
$ws.Range("B10").Font.Size = 10

$ws.Range("B11").Value = '渡部魁'
$ws.Range("B11").Font.Name = "Arial"
$ws.Range("B11").Font.Size = 10

$ws.Range("B12").Value = '崎谷航平'
$ws.Range("B12").Font.Name = "Arial"
$ws.Range("B12").Font.Size = 10

$ws.Range("B13").Value = '三神佳誠'
$ws.Range("B13").Font.Name = "Arial"
$ws.Range("B13").Font.Size = 10

$ws.Range("B14").Value = '氏家琉貴'
$ws.Range("B14").Font.Name = "Arial"
$ws.Range("B14").Font.Size = 10

$ws.Range("B15").Value = '羽賀尚生'
$ws.Range("B15").Font.Name = "Arial"
$ws.Range("B15").Font.Size = 10

$ws.Range("B16").Value = '足立耕平'
$ws.Range("B16").Font.Name = "Arial"
$ws.Range("B16").Font.Size = 10

$ws.Range("B17").Value = '遠藤隼人'
$ws.Range("B17").Font.Name = "Arial"
$ws.Range("B17").Font.Size = 10

$ws.Range("B18").Value = 'Ethan Virtudazo'
$ws.Range("B18").Font.Name = "Roboto"
$ws.Range("B18").Font.Size = 10

$ws.Range("B19").Value = '富澤天音'
$ws.Range("B19").Font.Name = "Arial"
$ws.Range("B19").Font.Size = 10

$ws.Range("B20").Value = 'みな'
$ws.Range("B20").Font.Name = "ArialMT"
$ws.Range("B20").Font.Size = 12

$ws.Range("B21").Value = '池田伊吹'
$ws.Range("B21").Font.Name = "Arial"
$ws.Range("B21").Font.Size = 10

$ws.Range("B22").Value = '神山修造'
$ws.Range("B22").Font.Name = "Arial"
$ws.Range("B22").Font.Size = 10

$ws.Range("B23").Value = '川田涼介'
$ws.Range("B23").Font.Name = "Arial"
$ws.Range("B23").Font.Size = 10

$ws.Range("B24").Value = '豊島亮'
$ws.Range("B24").Font.Name = "Arial"
$ws.Range("B24").Font.Size = 10

$ws.Range("B25").Value = '兒島大志郎'
$ws.Range("B25").Font.Name = "Arial"
$ws.Range("B25").Font.Size = 10

$ws.Range("B26").Value = '高野怜央'
$ws.Range("B26").Font.Name = "Arial"
$ws.Range("B26").Font.Size = 10

$ws.Range("B27").Value = '山口玲'
$ws.Range("B27").Font.Name = "Arial"
$ws.Range("B27").Font.Size = 10

$ws.Range("B28").Value = '志塚惇希'
$ws.Range("B28").Font.Name = "Arial"
$ws.Range("B28").Font.Size = 10

$ws.Range("B29").Value = '山口洸翔'
$ws.Range("B29").Font.Name = "Arial"
$ws.Range("B29").Font.Size = 10

$ws.Range("B30").Value = 'Owen Williamson'
$ws.Range("B30").Font.Name = "Arial"
$ws.Range("B30").Font.Size = 10

$ws.Range("B31").Value = '白岩詩佑介'
$ws.Range("B31").Font.Name = "Arial"
$ws.Range("B31").Font.Size = 10

$ws.Range("B32").Value = '石井海成'
$ws.Range("B32").Font.Name = "Arial"
$ws.Range("B32").Font.Size = 10

# View: zoom + active selection, matching the saved window state in the target workbook.
$ws.Application.ActiveWindow.Zoom = 88
$ws.Range("D29").Select() | Out-Null

# Page setup: paper size (A4=9 per OOXML paperSize code) + portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "edit complete"
